$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "30.303.93"
$ws.Range("E2").Value = "  -0.64%  "

# Row 3
$ws.Range("D3").Value = "1.883.61"
$ws.Range("E3").Value = "  -1.45%  "

# Row 4
$ws.Range("D4").Value = "0.9997"
$ws.Range("E4").Value = "  -0.03%  "

# Row 5
$ws.Range("D5").Value = "238.05"
$ws.Range("E5").Value = "  -0.55%  "

# Row 6
$ws.Range("E6").Value = "  -0.03%  "

# Row 7
$ws.Range("D7").Value = "0.4685"
$ws.Range("E7").Value = "  -1.92%  "

# Row 8
$ws.Range("D8").Value = "0.2835"

# Row 9
$ws.Range("D9").Value = "0.06604"
$ws.Range("E9").Value = "  -1.56%  "

# Row 10
$ws.Range("D10").Value = "20.72"
$ws.Range("E10").Value = "  +10.30%  "

# Row 11
$ws.Range("D11").Value = "0.07783"
$ws.Range("E11").Value = "  +1.28%  "

# Row 12
$ws.Range("D12").Value = "98.30"
$ws.Range("E12").Value = "  -3.48%  "

# Row 13
$ws.Range("D13").Value = "1.882.73"
$ws.Range("E13").Value = "  -1.69%  "

# Row 14
$ws.Range("D14").Value = "5.088"
$ws.Range("E14").Value = "  -2.05%  "

# Row 15
$ws.Range("D15").Value = "0.6780"
$ws.Range("E15").Value = "  +0.98%  "

# Row 16
$ws.Range("D16").Value = "285.57"
$ws.Range("E16").Value = "  +9.97%  "

# Row 17
$ws.Range("D17").Value = "30.303.04"
$ws.Range("E17").Value = "  -0.70%  "

# Row 18
$ws.Range("D18").Value = "1.000"
$ws.Range("E18").Value = "  +0.03%  "

# Row 19
$ws.Range("D19").Value = "12.66"
$ws.Range("E19").Value = "  -0.02%  "

# Row 20
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").Value = "5.423"
$ws.Range("E20").Value = "  +0.72%  "

# Row 21
$ws.Range("B21").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C21").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D21").Value = "2.132.36"
$ws.Range("E21").Value = "  -1.15%  "

# Row 22
$ws.Range("D22").Value = "0.000007295"
$ws.Range("E22").Value = "  -2.38%  "

# Row 23
$ws.Range("D23").Value = "0.9994"
$ws.Range("E23").Value = "  -0.05%  "

# Row 24
$ws.Range("D24").Value = "6.197"
$ws.Range("E24").Value = "  -1.27%  "

# Row 25
$ws.Range("D25").Value = "9.396"
$ws.Range("E25").Value = "  +0.43%  "

# Row 26
$ws.Range("D26").Value = "168.69"
$ws.Range("E26").Value = "  +0.64%  "

# Row 27
$ws.Range("D27").Value = "19.26"

# Row 28
$ws.Range("D28").Value = "1.997"
$ws.Range("E28").Value = "  -3.14%  "

# Row 29
$ws.Range("E29").Value = "  -0.78%  "

# Row 30
$ws.Range("D30").Value = "0.09715"
$ws.Range("E30").Value = "  -2.82%  "

# Row 31
$ws.Range("D31").Value = "4.393"
$ws.Range("E31").Value = "  -8.64%  "

# Row 32
$ws.Range("D32").Value = "1.486"
$ws.Range("E32").Value = "  -1.65%  "

# Row 33
$ws.Range("D33").Value = "4.141"
$ws.Range("E33").Value = "  -2.73%  "

# Row 34
$ws.Range("D34").Value = "0.04687"
$ws.Range("E34").Value = "  -0.75%  "

# Row 35
$ws.Range("D35").Value = "0.7079"
$ws.Range("E35").Value = "  -2.61%  "

# Row 36
$ws.Range("D36").Value = "1.099"
$ws.Range("E36").Value = "  -0.65%  "

# Row 37
$ws.Range("D37").Value = "2.714"
$ws.Range("E37").Value = "  +0.01%  "

# Row 38
$ws.Range("D38").Value = "0.01876"
$ws.Range("E38").Value = "  -2.15%  "

# Row 39
$ws.Range("D39").Value = "6.628"
$ws.Range("E39").Value = "  +5.91%  "

# Row 40
$ws.Range("D40").Value = "2.523"
$ws.Range("E40").Value = "  -3.84%  "

# Row 41
$ws.Range("D41").Value = "72.22"
$ws.Range("E41").Value = "  -3.64%  "

# Row 42
$ws.Range("D42").Value = "1.977"
$ws.Range("E42").Value = "  +0.34%  "

# Row 43
$ws.Range("D43").Value = "0.8658"
$ws.Range("E43").Value = "  +0.45%  "

# Row 44
$ws.Range("E44").Value = "  +0.00%  "

# Row 45
$ws.Range("D45").Value = "103.20"
$ws.Range("E45").Value = "  -2.26%  "

# Row 46
$ws.Range("D46").Value = "0.4191"
$ws.Range("E46").Value = "  -1.40%  "

# Row 47
$ws.Range("D47").Value = "985.51"
$ws.Range("E47").Value = "  +7.22%  "

# Row 48
$ws.Range("D48").Value = "7.316"
$ws.Range("E48").Value = "  -0.95%  "

# Row 49
$ws.Range("D49").Value = "9.207"
$ws.Range("E49").Value = "  +5.30%  "

# Row 50
$ws.Range("D50").Value = "33.96"
$ws.Range("E50").Value = "  -2.20%  "

# Row 51
$ws.Range("D51").Value = "0.1148"
$ws.Range("E51").Value = "  -4.45%  "
